# Switch back to 7% discount rate, adapt tax credits to reflect different
# discount rates for different technologies.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: production_tax_credits -------------------------------------
$ws1 = $wb.Worksheets.Item("production_tax_credits")

# Row 2 (nuclear_existing): discount rate 0.08 -> 0.07, and the CRF formula
# in F2 switches from an absolute reference ($E$2) to a relative one (E2).
$ws1.Range("E2").Value = 0.07
$ws1.Range("F2").Formula = "=E2/(1-(1/(1+E2)^C2))"

# Rows 3-18: each row now carries its own discount rate in column E (copied
# from the old shared $E$2), and the F-column CRF formula now references the
# row's own E cell (relative) instead of the old absolute $E$2.
$rows07 = 3,4,5,6,7,8,9,10,11,12,13,14
foreach ($r in $rows07) {
    $ws1.Range("E$r").Value = 0.07
    $ws1.Range("F$r").Formula = "=E$r/(1-(1/(1+E$r)^C$r))"
}

$rows10 = 15,16,17,18
foreach ($r in $rows10) {
    $ws1.Range("E$r").Value = 0.1
    $ws1.Range("F$r").Formula = "=E$r/(1-(1/(1+E$r)^C$r))"
}

# Selection ends up on I18 after editing the last row.
$ws1.Range("I18").Select()

# --- Sheet 2: investment_tax_credits --------------------------------------
$ws2 = $wb.Worksheets.Item("investment_tax_credits")
$ws2.Range("E2").Value = 0.07
